$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update preparer name to uppercase (B2:B41 and E2:E41), fixing style/font as a
# side effect of Excel baking the theme color to an explicit RGB black.
for ($r = 2; $r -le 41; $r++) {
    $ws.Cells.Item($r, 2).Value = "J.PLAGGENBERG"
    $ws.Cells.Item($r, 2).Font.Color = 0
    $ws.Cells.Item($r, 5).Value = "J.PLAGGENBERG"
    $ws.Cells.Item($r, 5).Font.Color = 0
}

$ws.Range("K11").Select()
